$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 - VENDA 19 (20/10)
$ws.Cells.Item(23, 2).Value = "CESAR AUGUSTO "
$ws.Cells.Item(23, 3).Value = "b3586279d9d3129597b2778d61367179"
$ws.Cells.Item(23, 4).Value = Get-Date -Year 2022 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(23, 5).Value = 365
$ws.Cells.Item(23, 6).Value = "-"
$ws.Cells.Item(23, 7).Value = "VENDA 19 (20/10)"

# Row 24 - VENDA 20 (20/10)
$ws.Cells.Item(24, 2).Value = "FERNANDO BRITO"
$ws.Cells.Item(24, 3).Value = "7687deae1489a0478188e4675601a77b"
$ws.Cells.Item(24, 4).Value = Get-Date -Year 2022 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(24, 5).Value = 365
$ws.Cells.Item(24, 6).Value = "-"
$ws.Cells.Item(24, 7).Value = "VENDA 20 (20/10)"

$ws.Range("D23:D24").NumberFormat = "yyyy-mm-dd"
